# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other header cells (e.g. column G's "sum" header),
# and a corresponding 0 value in the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save"
$ws.Range("H1").Value = "Save"

# Copy the formatting (bold font, borders, centered alignment) from the
# neighboring header cell G1 so H1 matches the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell H2 = 0
$ws.Range("H2").Value = 0
